$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free inline pattern: force text format so numeric-looking
# strings (prices like "51.749.75", "1.00") and percent strings
# ("  +5.22%  ") are stored as literal text, matching the source data,
# then restore the cell style to Normal so no stray formatting is left
# behind (matches the originally unstyled data cells).

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '51.749.75'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +5.22%  '
$ws.Range('E2').Style = "Normal"

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.761.50'
$ws.Range('D3').Style = "Normal"

# Row 4
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('E4').Style = "Normal"

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '116.55'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +4.27%  '
$ws.Range('E5').Style = "Normal"

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '332.77'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +2.95%  '
$ws.Range('E6').Style = "Normal"

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.539'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +2.34%  '
$ws.Range('E7').Style = "Normal"

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E8').Style = "Normal"

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.577'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +6.55%  '
$ws.Range('E9').Style = "Normal"

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '41.81'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +5.11%  '
$ws.Range('E10').Style = "Normal"

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0859'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +6.11%  '
$ws.Range('E11').Style = "Normal"

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '20.19'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +2.21%  '
$ws.Range('E12').Style = "Normal"

# Row 13
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +2.06%  '
$ws.Range('E13').Style = "Normal"

# Row 14
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +5.50%  '
$ws.Range('E14').Style = "Normal"

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.192.85'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +5.48%  '
$ws.Range('E15').Style = "Normal"

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '2.774.75'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +5.13%  '
$ws.Range('E16').Style = "Normal"

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.889'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +3.67%  '
$ws.Range('E17').Style = "Normal"

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '51.702.72'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +5.26%  '
$ws.Range('E18').Style = "Normal"

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '3.22'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +5.95%  '
$ws.Range('E19').Style = "Normal"

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '13.50'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +4.43%  '
$ws.Range('E20').Style = "Normal"

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.87'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +2.88%  '
$ws.Range('E21').Style = "Normal"

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.0₃0974'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +3.39%  '
$ws.Range('E22').Style = "Normal"

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '278.32'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +3.36%  '
$ws.Range('E23').Style = "Normal"

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '69.62'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +1.54%  '
$ws.Range('E24').Style = "Normal"

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.67'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +5.40%  '
$ws.Range('E25').Style = "Normal"

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '26.78'
$ws.Range('D26').Style = "Normal"

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.00'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('E27').Style = "Normal"

# Row 28
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -1.54%  '
$ws.Range('E28').Style = "Normal"

# Row 29
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +0.07%  '
$ws.Range('E29').Style = "Normal"

# Row 30
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +1.83%  '
$ws.Range('E30').Style = "Normal"

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '35.05'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +0.48%  '
$ws.Range('E31').Style = "Normal"

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '50.01'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +0.89%  '
$ws.Range('E32').Style = "Normal"

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.57'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +1.84%  '
$ws.Range('E33').Style = "Normal"

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0826'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +1.73%  '
$ws.Range('E34').Style = "Normal"

# Row 35
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('E35').Style = "Normal"

# Row 36
$ws.Range('B36').Value = 'Celestia'
$ws.Range('C36').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '18.99'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -0.06%  '
$ws.Range('E36').Style = "Normal"

# Row 37
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '5.02'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +2.67%  '
$ws.Range('E37').Style = "Normal"

# Row 38
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +2.48%  '
$ws.Range('E38').Style = "Normal"

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.24'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +3.62%  '
$ws.Range('E39').Style = "Normal"

# Row 40
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +9.64%  '
$ws.Range('E40').Style = "Normal"

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '126.86'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -0.77%  '
$ws.Range('E41').Style = "Normal"

# Row 42
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '23.26'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +4.16%  '
$ws.Range('E42').Style = "Normal"

# Row 43
$ws.Range('B43').Value = 'Stellar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.114'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +3.43%  '
$ws.Range('E43').Style = "Normal"

# Row 44
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +7.68%  '
$ws.Range('E44').Style = "Normal"

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.44'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +13.15%  '
$ws.Range('E45').Style = "Normal"

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.088.97'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +1.69%  '
$ws.Range('E46').Style = "Normal"

# Row 47
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +3.70%  '
$ws.Range('E47').Style = "Normal"

# Row 48
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +5.43%  '
$ws.Range('E48').Style = "Normal"

# Row 49
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +6.65%  '
$ws.Range('E49').Style = "Normal"

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '9.03'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +1.45%  '
$ws.Range('E50').Style = "Normal"

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '59.96'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +1.88%  '
$ws.Range('E51').Style = "Normal"
